# Update files for Project
# The header row (ESPAÑOL / QUECHUA) is removed, all remaining rows shift up
# by one, and a handful of misaligned Spanish/Quechua cells are corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift all data rows up by one (this drops the "ESPAÑOL"/"QUECHUA"
#     header that used to live in row 1) -------------------------------
$srcRange = $ws.Range("A2:B311")
$vals = $srcRange.Value2
$ws.Range("A1:B310").Value = $vals

# Remove the now-duplicated last row (311) so the sheet ends at row 310.
$ws.Rows("311").Delete()

# --- Fix a handful of rows where the Spanish/Quechua pairing had been
#     shifted/misaligned in the source data --------------------------
$ws.Range("A37").Value = '¿Le duele la lengua?'
$ws.Range("B37").Value = '¿Tonkori nanan?'

$ws.Range("A38").Value = '¿Kalloyki nanan?'
$ws.Range("B38").Value = '¿Arde, quema?'

$ws.Range("A41").Value = '¿En descanso o en actividad?'
$ws.Range("B41").Value = '¿Arde, quema?'

$ws.Range("A42").Value = '¿Samanki hina o ruwanki hina?'
$ws.Range("B42").Value = '¿Le duele el pulmón?'

$ws.Range("A224").Value = 'Abra la boca'
$ws.Range("B224").Value = 'Cierre los ojos'

$ws.Range("A225").Value = 'Şawiykita wiskay ó (chîrmiy)'
$ws.Range("B225").Value = 'K''alloykita urkûmuy'

# --- Cosmetic sheet view / column width changes -----------------------
$ws.Columns("A").ColumnWidth = 68

$ws.Range("E8").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 298
$excel.ActiveWindow.ScrollColumn = 1
